$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-7: Klassifizierung "Feature" -> "Leistung"
$ws.Range("C3").Value = "Leistung"
$ws.Range("C4").Value = "Leistung"
$ws.Range("C5").Value = "Leistung"
$ws.Range("C6").Value = "Leistung"
$ws.Range("C7").Value = "Leistung"

# Row 9: fill in new "Hauptziel" goal entry
$ws.Range("D9:E9").WrapText = $true
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "Hauptziel"
$ws.Range("D9").Value = "Funktionierendes Program"
$ws.Range("E9").Value = "Es soll ein vollständig funkionierendes Program vorliegen"
$ws.Range("F9").Value = "Muss"

# Update active selection
$ws.Range("M6").Select()
